# Video-Indexing.xlsx - "updated main GSC export data"
#
# Appends 7 more days of GSC video-indexing data (2025-11-21 .. 2025-11-27)
# to the "Chart" sheet, and rolls the "Videos" failure count on the
# "Table" sheet from 23 to 24 to stay consistent with the newest day
# (2025-11-27) which failed validation ("Video isn't on a watch page").

$wb = $excel.ActiveWorkbook

# --- helper: write a literal (non-date-autodetected) text value -----------
# Assigning a date-shaped string straight to .Value lets Excel's
# autodetection turn it into a real date serial + date number format, which
# is not what this sheet uses (the Date column is plain text). Routing the
# literal through a formula and then baking it back down to a static value
# with a self Copy/PasteSpecial(values) keeps the existing cell style (s=0)
# and produces a plain shared-string text cell instead.
function Set-LiteralText($cell, [string]$text) {
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# 1) "Chart" sheet: append rows 49-55 (2025-11-21 .. 2025-11-27)
# ---------------------------------------------------------------------
$chart = $wb.Worksheets.Item("Chart")

$newRows = @(
    @{ Date = "2025-11-21"; NoVideo = 24; VideoIdx = 1; Impr = 0 },
    @{ Date = "2025-11-22"; NoVideo = 24; VideoIdx = 1; Impr = 0 },
    @{ Date = "2025-11-23"; NoVideo = 24; VideoIdx = 1; Impr = 0 },
    @{ Date = "2025-11-24"; NoVideo = 24; VideoIdx = 1; Impr = 0 },
    @{ Date = "2025-11-25"; NoVideo = 24; VideoIdx = 1; Impr = 0 },
    @{ Date = "2025-11-26"; NoVideo = 24; VideoIdx = 1; Impr = 0 },
    @{ Date = "2025-11-27"; NoVideo = 24; VideoIdx = 1; Impr = $null }
)

$lastRow = 48
$r = $lastRow + 1
foreach ($row in $newRows) {
    # Carry over formatting (style 0, General) from the previous last row
    # for each of the four columns so no new cell styles get created.
    $chart.Cells.Item($lastRow, 1).Copy($chart.Cells.Item($r, 1))
    $chart.Cells.Item($lastRow, 2).Copy($chart.Cells.Item($r, 2))
    $chart.Cells.Item($lastRow, 3).Copy($chart.Cells.Item($r, 3))
    $chart.Cells.Item($lastRow, 4).Copy($chart.Cells.Item($r, 4))

    Set-LiteralText $chart.Cells.Item($r, 1) $row.Date
    $chart.Cells.Item($r, 2).Value = $row.NoVideo
    $chart.Cells.Item($r, 3).Value = $row.VideoIdx

    if ($null -eq $row.Impr) {
        # 2025-11-27 has no Impressions figure yet (blank cell), same as
        # the other not-yet-reported columns elsewhere on this sheet.
        $chart.Cells.Item($r, 4).ClearContents()
    } else {
        $chart.Cells.Item($r, 4).Value = $row.Impr
    }

    $r++
}

# ---------------------------------------------------------------------
# 2) "Table" sheet: Videos count 23 -> 24
# ---------------------------------------------------------------------
$table = $wb.Worksheets.Item("Table")
$table.Cells.Item(2, 3).Value = 24

Write-Host "Applied GSC export update: added rows 49-55 to Chart, Videos=24 on Table."
